$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.06279926495191
$ws.Range("C2").Value = 5.113490353836027
$ws.Range("D2").Value = 10.82885608225241
$ws.Range("F2").Value = 32.68360632238937
$ws.Range("G2").Value = 3.650477257462309
$ws.Range("I2").Value = 21.86596526605257
$ws.Range("J2").Value = 11.36552556472787
$ws.Range("K2").Value = 10.22771242244629
$ws.Range("N2").Value = 18.68885310334556
$ws.Range("O2").Value = 23.95311639885094
$ws.Range("B3").Value = 9.782878767430086
$ws.Range("C3").Value = 4.909417274066606
$ws.Range("D3").Value = 10.72922055032913
$ws.Range("F3").Value = 32.71302454322903
$ws.Range("G3").Value = 3.652215704730075
$ws.Range("I3").Value = 21.95090625029659
$ws.Range("J3").Value = 11.34658089783536
$ws.Range("K3").Value = 10.04100689154385
$ws.Range("N3").Value = 18.74582830612414
$ws.Range("O3").Value = 24.02439165531005
$ws.Range("B4").Value = 9.608475668004994
$ws.Range("C4").Value = 4.780336197466132
$ws.Range("D4").Value = 10.66985624810654
$ws.Range("F4").Value = 32.73888056914939
$ws.Range("G4").Value = 3.653340024894119
$ws.Range("I4").Value = 22.00695714460683
$ws.Range("J4").Value = 11.33720381199259
$ws.Range("K4").Value = 9.926196925185042
$ws.Range("N4").Value = 18.78246458191811
$ws.Range("O4").Value = 24.07294984911495
$ws.Range("B5").Value = 9.536878586286745
$ws.Range("C5").Value = 4.726862936773468
$ws.Range("D5").Value = 10.6461448224804
$ws.Range("F5").Value = 32.75137469412402
$ws.Range("G5").Value = 3.653812548495863
$ws.Range("I5").Value = 22.03077768942196
$ws.Range("J5").Value = 11.33395241064724
$ws.Range("K5").Value = 9.879428817461976
$ws.Range("N5").Value = 18.79781110428691
$ws.Range("O5").Value = 24.09394135343204
$ws.Range("B6").Value = 9.524961492605733
$ws.Range("C6").Value = 4.71793357376938
$ws.Range("D6").Value = 10.64223724134963
$ws.Range("F6").Value = 32.75356749037652
$ws.Range("G6").Value = 3.653891878932801
$ws.Range("I6").Value = 22.03479220805234
$ws.Range("J6").Value = 11.33344700764848
$ws.Range("K6").Value = 9.871665896017211
$ws.Range("N6").Value = 18.80038460068913
$ws.Range("O6").Value = 24.09749961168606
$ws.Range("B7").Value = 9.607512066026601
$ws.Range("C7").Value = 4.779618457100423
$ws.Range("D7").Value = 10.66953449258374
$ws.Range("F7").Value = 32.7390411464833
$ws.Range("G7").Value = 3.653346339334752
$ws.Range("I7").Value = 22.00727443214667
$ws.Range("J7").Value = 11.33715765184444
$ws.Range("K7").Value = 9.925566040169162
$ws.Range("N7").Value = 18.78266986069381
$ws.Range("O7").Value = 24.07322807718597
$ws.Range("B8").Value = 9.966871882755076
$ws.Range("C8").Value = 5.043954793697565
$ws.Range("D8").Value = 10.79414044628624
$ws.Range("F8").Value = 32.69213127354919
$ws.Range("G8").Value = 3.651064889903689
$ws.Range("I8").Value = 21.89444382048028
$ws.Range("J8").Value = 11.35852722182764
$ws.Range("K8").Value = 10.16340838278393
$ws.Range("N8").Value = 18.70815575058905
$ws.Range("O8").Value = 23.97669570907831
$ws.Range("B9").Value = 10.64700648483655
$ws.Range("C9").Value = 5.529319809781674
$ws.Range("D9").Value = 11.05176503879786
$ws.Range("F9").Value = 32.66203405698226
$ws.Range("G9").Value = 3.647040475204993
$ws.Range("I9").Value = 21.70412625260142
$ws.Range("J9").Value = 11.41818417709071
$ws.Range("K9").Value = 10.62583158267194
$ws.Range("N9").Value = 18.57509551920072
$ws.Range("O9").Value = 23.82553094098418
$ws.Range("B10").Value = 11.12622053703642
$ws.Range("C10").Value = 5.862349012896509
$ws.Range("D10").Value = 11.24760462573252
$ws.Range("F10").Value = 32.67768204094001
$ws.Range("G10").Value = 3.644354926806304
$ws.Range("I10").Value = 21.58319604744677
$ws.Range("J10").Value = 11.47261988836791
$ws.Range("K10").Value = 10.95975890334442
$ws.Range("N10").Value = 18.48521898965278
$ws.Range("O10").Value = 23.73783609148637
$ws.Range("B11").Value = 11.33876476992045
$ws.Range("C11").Value = 6.008155001147217
$ws.Range("D11").Value = 11.33778819698629
$ws.Range("F11").Value = 32.69298692661508
$ws.Range("G11").Value = 3.64319148387096
$ws.Range("I11").Value = 21.53229234860199
$ws.Range("J11").Value = 11.49963196159948
$ws.Range("K11").Value = 11.10973506914537
$ws.Range("N11").Value = 18.44602647883654
$ws.Range("O11").Value = 23.70303825612419
$ws.Range("B12").Value = 11.41839100667139
$ws.Range("C12").Value = 6.062509987869881
$ws.Range("D12").Value = 11.37206756898328
$ws.Range("F12").Value = 32.69995684890924
$ws.Range("G12").Value = 3.642759245548089
$ws.Range("I12").Value = 21.51360792373008
$ws.Range("J12").Value = 11.51017876175495
$ws.Range("K12").Value = 11.16619708719788
$ws.Range("N12").Value = 18.43142745902227
$ws.Range("O12").Value = 23.6905955918881
$ws.Range("B13").Value = 11.40128143231928
$ws.Range("C13").Value = 6.050842434870826
$ws.Range("D13").Value = 11.36467961984478
$ws.Range("F13").Value = 32.69840358161505
$ws.Range("G13").Value = 3.642851965868363
$ws.Range("I13").Value = 21.51760561725431
$ws.Range("J13").Value = 11.50789327210602
$ws.Range("K13").Value = 11.15405250752637
$ws.Range("N13").Value = 18.43456085846041
$ws.Range("O13").Value = 23.69324264595814
$ws.Range("B14").Value = 11.34533330442906
$ws.Range("C14").Value = 6.012644236064682
$ws.Range("D14").Value = 11.34060596401259
$ws.Range("F14").Value = 32.69353683202475
$ws.Range("G14").Value = 3.643155756605084
$ws.Range("I14").Value = 21.53074330554075
$ws.Range("J14").Value = 11.50049331462922
$ws.Range("K14").Value = 11.11438715620111
$ws.Range("N14").Value = 18.44482055814218
$ws.Range("O14").Value = 23.70199985954449
$ws.Range("B15").Value = 11.31094938828516
$ws.Range("C15").Value = 5.989133827415016
$ws.Range("D15").Value = 11.3258760796078
$ws.Range("F15").Value = 32.69070863126319
$ws.Range("G15").Value = 3.643342921109891
$ws.Range("I15").Value = 21.53886760495548
$ws.Range("J15").Value = 11.49600185885445
$ws.Range("K15").Value = 11.09004638310663
$ws.Range("N15").Value = 18.45113644739659
$ws.Range("O15").Value = 23.70745961512532
$ws.Range("B16").Value = 11.11221392619229
$ws.Range("C16").Value = 5.852702162821713
$ws.Range("D16").Value = 11.24173050306997
$ws.Range("F16").Value = 32.67684633817422
$ws.Range("G16").Value = 3.644432128723653
$ws.Range("I16").Value = 21.58660547075502
$ws.Range("J16").Value = 11.47089938175064
$ws.Range("K16").Value = 10.94991416895428
$ws.Range("N16").Value = 18.48781428229047
$ws.Range("O16").Value = 23.74021289879857
$ws.Range("B17").Value = 10.98884354225925
$ws.Range("C17").Value = 5.767517718785557
$ws.Range("D17").Value = 11.19037001086743
$ws.Range("F17").Value = 32.67043728180671
$ws.Range("G17").Value = 3.645115206500543
$ws.Range("I17").Value = 21.61694395944298
$ws.Range("J17").Value = 11.45607207485737
$ws.Range("K17").Value = 10.86341469187485
$ws.Range("N17").Value = 18.51074772992912
$ws.Range("O17").Value = 23.76161229794557
$ws.Range("B18").Value = 10.91737498322429
$ws.Range("C18").Value = 5.717988260381146
$ws.Range("D18").Value = 11.16093361516276
$ws.Range("F18").Value = 32.66752175641955
$ws.Range("G18").Value = 3.645513577732604
$ws.Range("I18").Value = 21.63478046641204
$ws.Range("J18").Value = 11.44775584164962
$ws.Range("K18").Value = 10.81348392177067
$ws.Range("N18").Value = 18.52409783838412
$ws.Range("O18").Value = 23.77440012763396
$ws.Range("B19").Value = 10.89309192629546
$ws.Range("C19").Value = 5.701128110098941
$ws.Range("D19").Value = 11.15098583097672
$ws.Range("F19").Value = 32.66666706678105
$ws.Range("G19").Value = 3.64564940233468
$ws.Range("I19").Value = 21.64088597282665
$ws.Range("J19").Value = 11.44497668669209
$ws.Range("K19").Value = 10.79654922414116
$ws.Range("N19").Value = 18.52864536539993
$ws.Range("O19").Value = 23.77881215592108
$ws.Range("B20").Value = 11.00202979837806
$ws.Range("C20").Value = 5.776641285910136
$ws.Range("D20").Value = 11.19582678424006
$ws.Range("F20").Value = 32.67103978040826
$ws.Range("G20").Value = 3.645041924534057
$ws.Range("I20").Value = 21.61367435581533
$ws.Range("J20").Value = 11.45762855992034
$ws.Range("K20").Value = 10.87264158423158
$ws.Range("N20").Value = 18.50828993628544
$ws.Range("O20").Value = 23.7592846589611
$ws.Range("B21").Value = 11.36179052322397
$ws.Range("C21").Value = 6.023887567032677
$ws.Range("D21").Value = 11.34767371260802
$ws.Range("F21").Value = 32.69493447261705
$ws.Range("G21").Value = 3.643066300088703
$ws.Range("I21").Value = 21.5268683816235
$ws.Range("J21").Value = 11.50265827921002
$ws.Range("K21").Value = 11.12604721708121
$ws.Range("N21").Value = 18.44180046737788
$ws.Range("O21").Value = 23.69940770474247
$ws.Range("B22").Value = 11.59187418532195
$ws.Range("C22").Value = 6.180458096496293
$ws.Range("D22").Value = 11.4476520402183
$ws.Range("F22").Value = 32.71739365312097
$ws.Range("G22").Value = 3.641823663519121
$ws.Range("I22").Value = 21.47358500891368
$ws.Range("J22").Value = 11.53393801689753
$ws.Range("K22").Value = 11.2897121956136
$ws.Range("N22").Value = 18.39975776526132
$ws.Range("O22").Value = 23.66455677254643
$ws.Range("B23").Value = 11.46955856410994
$ws.Range("C23").Value = 6.097364567313408
$ws.Range("D23").Value = 11.39423368629805
$ws.Range("F23").Value = 32.70478189146577
$ws.Range("G23").Value = 3.642482453215484
$ws.Range("I23").Value = 21.50170740739365
$ws.Range("J23").Value = 11.51707605696725
$ws.Range("K23").Value = 11.20255642673157
$ws.Range("N23").Value = 18.4220679032281
$ws.Range("O23").Value = 23.6827649752098
$ws.Range("B24").Value = 10.996069972132
$ws.Range("C24").Value = 5.772518249247784
$ws.Range("D24").Value = 11.19335949000513
$ws.Range("F24").Value = 32.67076499478655
$ws.Range("G24").Value = 3.645075037675493
$ws.Range("I24").Value = 21.61515131481111
$ws.Range("J24").Value = 11.45692422414823
$ws.Range("K24").Value = 10.86847072838273
$ws.Range("N24").Value = 18.50940058959709
$ws.Range("O24").Value = 23.76033547383113
$ws.Range("B25").Value = 10.46623926721995
$ws.Range("C25").Value = 5.401934839147743
$ws.Range("D25").Value = 10.98080565515455
$ws.Range("F25").Value = 32.66354337060575
$ws.Range("G25").Value = 3.648081358543657
$ws.Range("I25").Value = 21.75229627664601
$ws.Range("J25").Value = 11.40016550914095
$ws.Range("K25").Value = 10.50151813955274
$ws.Range("N25").Value = 18.60970179855409
$ws.Range("O25").Value = 23.8623298062519
